$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.76699298620224
$ws.Range("B1").Value = 1.601200342178345
$ws.Range("C1").Value = 4.200932502746582
$ws.Range("D1").Value = 3.713629722595215
$ws.Range("E1").Value = 1.79625403881073
